{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  para.load(\"text\");\n}\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text.indexOf(\"Distance_Range_4_point_corre_function_average\") !== -1) {\n    para.insertText(\"  12  Distance_Range_4_point_corre_function_average 13 hbar_scale\", Word.InsertLocation.replace);\n    break;\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$rng = $d.Content\n$rng.Find.Text = \"Distance_Range_4_point_corre_function_average\"\n$rng.Find.Replacement.Text = \"Distance_Range_4_point_corre_function_average 13 hbar_scale\"\n$null = $rng.Find.Execute($rng.Find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $rng.Find.Replacement.Text, 2)\n"}
